$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-03-21 Thursday" "2024-03-22 Friday"

Replace-Text "75÷9=8, 3" "33÷2=16, 1"
Replace-Text "99÷9=11, 0" "79÷3=26, 1"
Replace-Text "95÷3=31, 2" "69÷4=17, 1"
Replace-Text "74÷7=10, 4" "87÷9=9, 6"
Replace-Text "83÷5=16, 3" "55÷2=27, 1"

Replace-Text "14÷6=2, 2" "53÷2=26, 1"
Replace-Text "92÷9=10, 2" "52÷3=17, 1"
Replace-Text "65÷6=10, 5" "68÷5=13, 3"
Replace-Text "53÷4=13, 1" "95÷9=10, 5"
Replace-Text "30÷2=15, 0" "89÷3=29, 2"

Replace-Text "71÷8=8, 7" "80÷2=40, 0"
Replace-Text "41÷7=5, 6" "61÷9=6, 7"
Replace-Text "21÷4=5, 1" "24÷3=8, 0"
Replace-Text "80÷5=16, 0" "36÷3=12, 0"
Replace-Text "73÷3=24, 1" "92÷9=10, 2"

Replace-Text "63÷4=15, 3" "21÷5=4, 1"
Replace-Text "96÷4=24, 0" "19÷3=6, 1"
Replace-Text "30÷7=4, 2" "68÷4=17, 0"
Replace-Text "97÷7=13, 6" "28÷6=4, 4"
Replace-Text "18÷7=2, 4" "61÷7=8, 5"

Replace-Text "32÷7=4, 4" "27÷7=3, 6"
Replace-Text "83÷3=27, 2" "22÷3=7, 1"
Replace-Text "22÷5=4, 2" "41÷2=20, 1"
Replace-Text "74÷6=12, 2" "77÷5=15, 2"
Replace-Text "55÷4=13, 3" "59÷9=6, 5"

Write-Output "Replacements complete"
